$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseRun")

$ws.Range("A18").Value = "TN2485417"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "11541021"
$ws.Range("C18").Value = "Personal Auto - Credit"
$ws.Range("D18").Value = "Restricted"
$ws.Range("E18").Value = "TC003"
